$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells: copy the existing header formatting (style index reused,
# matching B1:H1) onto the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF)
$iValues = @(8, 9, 7, 3, 9, 8, 2)
$jValues = @(8, 9, 8, 6, 9, 8, 2)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
